$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the last column (AD), which shifts the
# existing "ID" header from AD1 to AE1 and leaves a blank AD1 behind.
$ws.Columns("AD:AD").Insert() | Out-Null

# Give the newly inserted column its header text.
$ws.Range("AD1").Value = "Flag"

# Re-apply AutoFilter over the full, now-wider header row so the filter
# (and its underlying ref) covers the new column too.
$ws.AutoFilterMode = $false
$ws.Range("A1:AE1").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$AE`$1"
    }
}

# Match the author's final selection: reviewing the newly added column.
$ws.Range("AD1").Select() | Out-Null
